# Updates to datasets 3 for JPF
# Adds a new "Higher/Lower" column (J) to each of the 4 worksheets, and
# corrects a handful of recomputed correlation values on row 9 of the
# "all_tools" and "infer" sheets.

$wb = $excel.ActiveWorkbook

# Same Higher/Lower pattern applies identically to every sheet.
# Rows 6, 7 and 8 have no correlation data (and therefore no
# Higher/Lower indicator) in any sheet.
$higherLower = @{
    2  = "Higher"
    3  = "Lower"
    4  = "Lower"
    5  = "Higher"
    9  = "Lower"
    10 = "Higher"
    11 = "Lower"
    12 = "Lower"
    13 = "Higher"
    14 = "Lower"
    15 = "Lower"
    16 = "Lower"
    17 = "Lower"
    18 = "Higher"
    19 = "Lower"
}

$sheetNames = @("all_tools", "checker_framework", "typestate_checker", "infer")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # New header cell J1, matching the style used by the other header cells.
    $ws.Range("I1").Copy()
    $ws.Range("J1").PasteSpecial(-4122)
    $ws.Range("J1").Value = "Higher/Lower"

    foreach ($row in $higherLower.Keys) {
        $ws.Cells.Item($row, 10).Value = $higherLower[$row]
    }
}

# Corrected correlation figures on row 9 of the "all_tools" sheet.
$wsAll = $wb.Worksheets.Item("all_tools")
$wsAll.Range("D9").Value = 392
$wsAll.Range("F9").Value = -0.2296446241402438
$wsAll.Range("G9").Value = 0.001381085952098848
$wsAll.Range("H9").Value = -0.3250508508876161
$wsAll.Range("I9").Value = 0.0009677863595653327

# Corrected correlation figures on row 9 of the "infer" sheet.
$wsInfer = $wb.Worksheets.Item("infer")
$wsInfer.Range("C9").Value = 13
$wsInfer.Range("D9").Value = 13
$wsInfer.Range("F9").Value = -0.1325530043077417
$wsInfer.Range("G9").Value = 0.1086826442074059
$wsInfer.Range("H9").Value = -0.1612223880273475
$wsInfer.Range("I9").Value = 0.1090548020620709
